$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new Price (D) [or $null to skip], new Volume(1h) % (E) [or $null to skip]
# Values are prefixed with a leading apostrophe so Excel always stores them
# as literal text (matching the original inline-string cell type), even
# when the digits would otherwise be auto-recognised as a number. The
# style is then reset to "Normal" so no stray quote-prefix/number format
# is left applied to the cell.
$updates = @(
    @{ Row = 2;  D = "38.780.68";  E = "  +1.49%  " },
    @{ Row = 3;  D = "2.102.41";   E = "  +0.48%  " },
    @{ Row = 4;  D = $null;        E = "  -0.08%  " },
    @{ Row = 5;  D = "229.35";     E = "  +0.40%  " },
    @{ Row = 6;  D = $null;        E = "  +1.02%  " },
    @{ Row = 7;  D = "61.93";      E = "  +1.65%  " },
    @{ Row = 8;  D = $null;        E = "  -0.05%  " },
    @{ Row = 9;  D = $null;        E = "  +2.02%  " },
    @{ Row = 10; D = $null;        E = "  -0.87%  " },
    @{ Row = 11; D = $null;        E = "  +0.23%  " },
    @{ Row = 12; D = "15.66";      E = "  +6.52%  " },
    @{ Row = 13; D = "2.411.59";   E = "  +0.40%  " },
    @{ Row = 14; D = "21.97";      E = "  -1.45%  " },
    @{ Row = 16; D = $null;        E = "  +0.36%  " },
    @{ Row = 17; D = "2.085.45";   E = "  -0.47%  " },
    @{ Row = 18; D = "38.773.63";  E = "  +1.69%  " },
    @{ Row = 19; D = "72.03";      E = "  +2.58%  " },
    @{ Row = 20; D = "6.06";       E = "  +0.64%  " },
    @{ Row = 21; D = "0.0₃0840";   E = "  +0.37%  " },
    @{ Row = 22; D = "227.81";     E = "  +1.68%  " },
    @{ Row = 23; D = $null;        E = "  +0.05%  " },
    @{ Row = 24; D = $null;        E = "  -0.87%  " },
    @{ Row = 25; D = $null;        E = "  +0.82%  " },
    @{ Row = 26; D = "171.92";     E = "  +1.23%  " },
    @{ Row = 27; D = "9.57";       E = "  +1.10%  " },
    @{ Row = 28; D = $null;        E = "  +5.70%  " },
    @{ Row = 29; D = $null;        E = "  +3.63%  " },
    @{ Row = 30; D = "19.34";      E = $null },
    @{ Row = 31; D = "2.47";       E = "  +3.67%  " },
    @{ Row = 32; D = $null;        E = "  +1.28%  " },
    @{ Row = 33; D = $null;        E = "  +2.36%  " },
    @{ Row = 34; D = $null;        E = "  +1.52%  " },
    @{ Row = 35; D = "0.0622";     E = "  +2.79%  " },
    @{ Row = 36; D = $null;        E = "  +2.83%  " },
    @{ Row = 37; D = "2.43";       E = "  +1.66%  " },
    @{ Row = 38; D = $null;        E = "  +1.30%  " },
    @{ Row = 39; D = "1.00";       E = "  +0.01%  " },
    @{ Row = 40; D = "18.29";      E = "  +1.06%  " },
    @{ Row = 41; D = "0.0229";     E = "  +4.48%  " },
    @{ Row = 42; D = "101.56";     E = "  +1.48%  " },
    @{ Row = 43; D = "1.534.15";   E = "  -1.20%  " },
    @{ Row = 44; D = $null;        E = "  -1.01%  " },
    @{ Row = 45; D = $null;        E = "  +4.08%  " },
    @{ Row = 46; D = $null;        E = "  -0.33%  " },
    @{ Row = 47; D = "1.14";       E = "  +1.62%  " },
    @{ Row = 48; D = "4.11";       E = "  -1.52%  " },
    @{ Row = 49; D = "1.05";       E = "  +2.00%  " },
    @{ Row = 50; D = $null;        E = "  -0.80%  " },
    @{ Row = 51; D = "2.294.08";   E = "  +0.24%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($r, 4)
        $cell.Value = "'" + $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $cell = $ws.Cells.Item($r, 5)
        $cell.Value = "'" + $u.E
        $cell.Style = "Normal"
    }
}
